$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-number-looking text cells to stay as Text so Excel
# does not silently coerce them into floating point numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = "51.759.65"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "3.101.86"
$ws.Range("E3").Value = "  +3.83%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "387.79"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("D6").Value = "103.64"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D11").Value = "0.137"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "0.0863"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "3.588.77"
$ws.Range("E13").Value = "  +3.86%  "
$ws.Range("D14").Value = "18.73"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "3.098.87"
$ws.Range("E16").Value = "  +3.36%  "
$ws.Range("D17").Value = "0.982"
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("D18").Value = "10.69"
$ws.Range("E18").Value = "  -4.65%  "
$ws.Range("D19").Value = "51.916.53"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("D20").Value = "3.18"
$ws.Range("E20").Value = "  +2.06%  "
$ws.Range("D21").Value = "12.53"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").Value = "70.20"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "268.93"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").Value = "3.13"
$ws.Range("E25").Value = "  -3.15%  "
$ws.Range("D26").Value = "8.19"
$ws.Range("E26").Value = "  +4.46%  "
$ws.Range("D27").Value = "27.13"
$ws.Range("E27").Value = "  +4.04%  "
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("D29").Value = "7.24"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "10.38"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("D33").Value = "35.64"
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("E34").Value = "  +2.88%  "
$ws.Range("D35").Value = "50.37"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("D39").Value = "0.295"
$ws.Range("E39").Value = "  +8.62%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "17.06"
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "1.89"
$ws.Range("E41").Value = "  +2.11%  "
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("D43").Value = "127.61"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("D46").Value = "22.22"
$ws.Range("E46").Value = "  +3.65%  "
$ws.Range("E47").Value = "  +4.48%  "
$ws.Range("D48").Value = "2.08"
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("D49").Value = "2.053.70"
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("D50").Value = "3.409.84"
$ws.Range("E50").Value = "  +3.91%  "
$ws.Range("E51").Value = "  +7.24%  "
